$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a value to be written as TEXT (avoids Excel auto-converting
# numeric-looking strings like "582.39" or "139.00" into float cells), while
# leaving the cell's style/format untouched afterwards.
function Set-TextValue($addr, $value) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "62.073.71"
$ws.Range("E2").Value = "  +2.42%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.410.77"
$ws.Range("E3").Value = "  +1.87%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.05%  "

# Row 5 - BNB
Set-TextValue "D5" "582.39"
$ws.Range("E5").Value = "  +1.61%  "

# Row 6 - Solana
Set-TextValue "D6" "139.00"
$ws.Range("E6").Value = "  +4.70%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.05%  "

# Row 8 - LidoStakedEther
Set-TextValue "D8" "3.409.24"
$ws.Range("E8").Value = "  +1.87%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +0.03%  "

# Row 10 - Toncoin
$ws.Range("E10").Value = "  -0.71%  "

# Row 11 - Dogecoin
Set-TextValue "D11" "0.128"
$ws.Range("E11").Value = "  +5.51%  "

# Row 12 - Cardano
$ws.Range("E12").Value = "  +2.06%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "3.992.05"
$ws.Range("E13").Value = "  +2.01%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  +1.88%  "

# Row 15 - ShibaInu
Set-TextValue "D15" "0.0000179"
$ws.Range("E15").Value = "  +3.89%  "

# Row 16 - WrappedEther
Set-TextValue "D16" "3.432.62"
$ws.Range("E16").Value = "  +2.59%  "

# Row 17 - Avalanche
Set-TextValue "D17" "25.62"
$ws.Range("E17").Value = "  +2.61%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "62.161.68"
$ws.Range("E18").Value = "  +2.39%  "

# Row 19 - Chainlink
Set-TextValue "D19" "14.22"
$ws.Range("E19").Value = "  +2.17%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  +3.58%  "

# Row 21 - Polkadot
Set-TextValue "D21" "5.84"
$ws.Range("E21").Value = "  +1.29%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "396.82"
$ws.Range("E22").Value = "  +5.91%  "

# Row 23 - Polygon
Set-TextValue "D23" "0.568"
$ws.Range("E23").Value = "  +1.08%  "

# Row 24 - PEPE
Set-TextValue "D24" "0.0000132"
$ws.Range("E24").Value = "  +13.93%  "

# Row 25 - WrappedeETH
Set-TextValue "D25" "3.552.27"
$ws.Range("E25").Value = "  +2.06%  "

# Row 26 - Dai
Set-TextValue "D26" "0.998"
$ws.Range("E26").Value = "  -0.20%  "

# Row 27 - Litecoin
$ws.Range("E27").Value = "  +2.27%  "

# Row 28 - Fetch.AI
Set-TextValue "D28" "1.67"
$ws.Range("E28").Value = "  +0.08%  "

# Row 29 - RenderToken
Set-TextValue "D29" "7.73"
$ws.Range("E29").Value = "  +0.96%  "

# Row 30 - Binance-PegBSC-USD
Set-TextValue "D30" "0.998"
$ws.Range("E30").Value = "  -0.03%  "

# Row 31 - Kaspa
$ws.Range("E31").Value = "  +4.63%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +2.41%  "

# Row 33 - PancakeSwap
$ws.Range("E33").Value = "  +1.61%  "

# Row 34 - USDe
$ws.Range("E34").Value = "  +0.01%  "

# Row 35 - EthereumClassic
Set-TextValue "D35" "23.62"
$ws.Range("E35").Value = "  +2.12%  "

# Row 36 - RenzoRestakedETH
Set-TextValue "D36" "3.441.12"
$ws.Range("E36").Value = "  +1.89%  "

# Row 37 - NEARProtocol
Set-TextValue "D37" "5.44"
$ws.Range("E37").Value = "  -0.35%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "1.60"
$ws.Range("E38").Value = "  +3.65%  "

# Row 39 - Aptos
$ws.Range("E39").Value = "  +0.14%  "

# Row 40 - Monero
Set-TextValue "D40" "164.88"
$ws.Range("E40").Value = "  +1.53%  "

# Row 41 - Hedera
$ws.Range("E41").Value = "  +1.72%  "

# Row 42 - Stacks
$ws.Range("E42").Value = "  +12.85%  "

# Row 43 - ONDO
$ws.Range("E43").Value = "  +4.63%  "

# Row 44 - Mantle
$ws.Range("E44").Value = "  +5.36%  "

# Row 45 - FirstDigitalUSD
$ws.Range("E45").Value = "  +0.11%  "

# Row 46 - Filecoin
$ws.Range("E46").Value = "  +1.44%  "

# Rows 47 & 48 - swap OKB and EnergySwap entries
# (row 47 becomes EnergySwap, row 48 becomes OKB)
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D47" "25.11"
$ws.Range("E47").Value = "  +7.06%  "

$ws.Range("B48").Value = "OKB"
$ws.Range("C48").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D48" "41.79"
$ws.Range("E48").Value = "  +0.91%  "

# Row 49 - Cosmos
$ws.Range("E49").Value = "  +0.09%  "

# Row 50 - InjectiveProtocol
Set-TextValue "D50" "23.45"
$ws.Range("E50").Value = "  +2.44%  "

# Row 51 - Maker
Set-TextValue "D51" "2.360.85"
$ws.Range("E51").Value = "  +8.18%  "
